$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Requisitos" entry that used to reference LOM3206 (row 23,
# shared-string index 35) so it now references the new prerequisite
# LOM3263 instead.
$newReq = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Cells.Item(23, 2).Value = $newReq
$ws.Cells.Item(23, 3).Value = $newReq

# Remove the now-obsolete second requirement row (previously row 24,
# which held the LOM3221 requirement) entirely.
$ws.Rows.Item(24).Delete()
